# Reorder the monthly rows so that each calendar year's block starts with
# October, November, December, followed by January through September
# (e.g. 2014-01 moves from row 2 to row 5; 2014-10 moves from row 11 to row 2).
# Column A holds the month label, columns B:D hold the three price-index
# series. Column B is genuinely blank ("") for 2014/2015 rows both before
# and after the reorder, so it is left untouched for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, "2014-10", $null, 84.5534, 98.6481),
    @(3, "2014-11", $null, 82.5959, 98.1964),
    @(4, "2014-12", $null, 80.0243, 97.7389),
    @(5, "2014-01", $null, 99.8135, 98.70269999999999),
    @(6, "2014-02", $null, 97.64149999999999, 98.6157),
    @(7, "2014-03", $null, 95.03830000000001, 98.545),
    @(8, "2014-04", $null, 94.7749, 98.1313),
    @(9, "2014-05", $null, 94.27160000000001, 98.28700000000001),
    @(10, "2014-06", $null, 93.1816, 98.9303),
    @(11, "2014-07", $null, 91.3809, 99.6584),
    @(12, "2014-08", $null, 89.0249, 99.4431),
    @(13, "2014-09", $null, 86.6992, 99.1138),
    @(14, "2015-10", $null, 80.90000000000001, 91.3),
    @(15, "2015-11", $null, 81.46469999999999, 91.1853),
    @(16, "2015-12", $null, 81.7946, 91.7153),
    @(17, "2015-01", $null, 78.5295, 97.49769999999999),
    @(18, "2015-02", $null, 78.0046, 96.8475),
    @(19, "2015-03", $null, 77.5936, 96.62350000000001),
    @(20, "2015-04", $null, 76.0714, 95.7804),
    @(21, "2015-05", $null, 76.47190000000001, 94.83620000000001),
    @(22, "2015-06", $null, 79.0244, 94.3095),
    @(23, "2015-07", $null, 79.10290000000001, 93.3099),
    @(24, "2015-08", $null, 78.8051, 93.0902),
    @(25, "2015-09", $null, 79.3369, 92.5502),
    @(26, "2016-10", 94.5, 105, 101),
    @(27, "2016-11", 94.5, 110.4, 106.1),
    @(28, "2016-12", 97.2, 118.6, 107.7),
    @(29, "2016-01", 97.25839999999999, 80.79040000000001, 90.9314),
    @(30, "2016-02", 97.92749999999999, 81.2205, 91.5702),
    @(31, "2016-03", 97.95650000000001, 86.66379999999999, 92.39449999999999),
    @(32, "2016-04", 98.0549, 92.4538, 94.2205),
    @(33, "2016-05", 99.09999999999999, 95.90000000000001, 95.09999999999999),
    @(34, "2016-06", 99, 94.2, 95.8),
    @(35, "2016-07", 98.90000000000001, 96.3, 97.90000000000001),
    @(36, "2016-08", 96.3, 99.59999999999999, 97.5),
    @(37, "2016-09", 93.5, 103.5, 98),
    @(38, "2017-10", 105.9, 111.5, 106.3),
    @(39, "2017-11", 105.9, 107.7, 100.9),
    @(40, "2017-12", 102.9, 105.1, 99),
    @(41, "2017-01", 100, 124.2, 110.1),
    @(42, "2017-02", 100, 128.3, 110.2),
    @(43, "2017-03", 100, 129.8, 109.7),
    @(44, "2017-04", 100, 122.3, 109.2),
    @(45, "2017-05", 100, 114.1, 109.2),
    @(46, "2017-06", 100, 112.1, 108.2),
    @(47, "2017-07", 100, 112.6, 106.8),
    @(48, "2017-08", 102.9, 114.4, 107.2),
    @(49, "2017-09", 105.9, 113.9, 107.5)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    if ($r[2] -ne $null) {
        $ws.Cells.Item($rowNum, 2).Value = $r[2]
    }
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}
